$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.723.33'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '1.696.01'
$ws.Range('E3').Value = '  +0.39%  '
$ws.Range('D4').Value = "'1.000"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'315.45"
$ws.Range('E5').Value = '  +1.28%  '
$ws.Range('D6').Value = "'1.000"
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('D7').Value = "'0.3942"
$ws.Range('E7').Value = '  +0.56%  '
$ws.Range('D8').Value = "'0.4051"
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').Value = "'1.490"
$ws.Range('E9').Value = '  +1.35%  '
$ws.Range('D10').Value = "'0.9998"
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('D11').Value = "'53.10"
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('D12').Value = "'0.08871"
$ws.Range('E12').Value = '  +1.62%  '
$ws.Range('D13').Value = "'7.236"
$ws.Range('E13').Value = '  -0.72%  '
$ws.Range('D14').Value = "'23.64"
$ws.Range('E14').Value = '  +2.87%  '
$ws.Range('D15').Value = "'8.080"
$ws.Range('E15').Value = '  +8.87%  '
$ws.Range('D16').Value = "'0.00001322"
$ws.Range('E16').Value = '  +0.41%  '
$ws.Range('D17').Value = '1.694.75'
$ws.Range('E17').Value = '  +0.78%  '
$ws.Range('D18').Value = "'99.86"
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = "'0.06999"
$ws.Range('E19').Value = '  -0.27%  '
$ws.Range('D20').Value = "'19.63"
$ws.Range('E20').Value = '  +1.68%  '
$ws.Range('D21').Value = "'7.063"
$ws.Range('E21').Value = '  +5.41%  '
$ws.Range('D22').Value = "'0.9997"
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').Value = "'14.35"
$ws.Range('E23').Value = '  +1.40%  '
$ws.Range('D24').Value = '24.724.59'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').Value = "'3.251"
$ws.Range('E25').Value = '  +9.59%  '
$ws.Range('D26').Value = "'2.350"
$ws.Range('E26').Value = '  +1.59%  '
$ws.Range('D27').Value = "'22.69"
$ws.Range('E27').Value = '  +1.54%  '
$ws.Range('E28').Value = '  +2.54%  '
$ws.Range('D29').Value = "'136.04"
$ws.Range('E29').Value = '  +1.97%  '
$ws.Range('D30').Value = "'5.188"
$ws.Range('E30').Value = '  +1.77%  '
$ws.Range('D31').Value = "'7.452"
$ws.Range('E31').Value = '  +1.88%  '
$ws.Range('D32').Value = '1.884.66'
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = "'1.064"
$ws.Range('E33').Value = '  -0.07%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = "'0.08597"
$ws.Range('E34').Value = '  -1.05%  '
$ws.Range('D35').Value = "'7.151"
$ws.Range('E35').Value = '  -2.85%  '
$ws.Range('D36').Value = "'11.60"
$ws.Range('E36').Value = '  +5.79%  '
$ws.Range('D37').Value = "'0.2748"
$ws.Range('E37').Value = '  +2.05%  '
$ws.Range('E38').Value = '  -0.81%  '
$ws.Range('D39').Value = "'14.51"
$ws.Range('E39').Value = '  -0.35%  '
$ws.Range('D40').Value = "'0.09226"
$ws.Range('E40').Value = '  +3.47%  '
$ws.Range('D41').Value = "'0.02737"
$ws.Range('E41').Value = '  +0.42%  '
$ws.Range('E42').Value = '  +1.55%  '
$ws.Range('D43').Value = "'0.7674"
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('D44').Value = "'16.06"
$ws.Range('E44').Value = '  +4.60%  '
$ws.Range('D45').Value = "'0.7201"
$ws.Range('E45').Value = '  +1.20%  '
$ws.Range('D46').Value = "'2.575"
$ws.Range('E46').Value = '  +5.95%  '
$ws.Range('D47').Value = "'4.207"
$ws.Range('E47').Value = '  +1.81%  '
$ws.Range('D48').Value = "'0.9998"
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = "'1.329"
$ws.Range('D50').Value = "'139.48"
$ws.Range('E50').Value = '  -0.12%  '
$ws.Range('D51').Value = "'0.07998"
$ws.Range('E51').Value = '  +0.99%  '
